$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 210
$ws.Cells.Item(6, 6).Value = 17
$ws.Cells.Item(7, 6).Value = 56
$ws.Cells.Item(9, 6).Value = 371
$ws.Cells.Item(10, 6).Value = 4498
$ws.Cells.Item(11, 6).Value = 4498
$ws.Cells.Item(12, 6).Value = 126
$ws.Cells.Item(14, 6).Value = 1068
$ws.Cells.Item(16, 6).Value = 3936
$ws.Cells.Item(19, 6).Value = 37
$ws.Cells.Item(20, 6).Value = 157
$ws.Cells.Item(21, 6).Value = 3308
$ws.Cells.Item(22, 6).Value = 2
$ws.Cells.Item(23, 6).Value = 3
$ws.Cells.Item(24, 6).Value = 12
$ws.Cells.Item(25, 6).Value = 2796
$ws.Cells.Item(26, 6).Value = 105
$ws.Cells.Item(28, 6).Value = 5
$ws.Cells.Item(29, 6).Value = 124
$ws.Cells.Item(31, 6).Value = 155
$ws.Cells.Item(32, 6).Value = 66
$ws.Cells.Item(33, 6).Value = 38
$ws.Cells.Item(36, 6).Value = 123
$ws.Cells.Item(37, 6).Value = 5085
$ws.Cells.Item(38, 6).Value = 679
$ws.Cells.Item(41, 6).Value = 952
$ws.Cells.Item(42, 6).Value = 9
$ws.Cells.Item(43, 6).Value = 989
$ws.Cells.Item(44, 6).Value = 385
$ws.Cells.Item(46, 6).Value = 1876
$ws.Cells.Item(48, 6).Value = 48
$ws.Cells.Item(49, 6).Value = 669
$ws.Cells.Item(50, 6).Value = 797

# Sheet: 演出 (Sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 75
$ws.Cells.Item(7, 6).Value = 36
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(21, 6).Value = 703

# Sheet: 全部类型 (Sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(7, 6).Value = 210
$ws.Cells.Item(8, 6).Value = 17
$ws.Cells.Item(9, 6).Value = 75
$ws.Cells.Item(11, 6).Value = 371
$ws.Cells.Item(12, 6).Value = 4498
$ws.Cells.Item(13, 6).Value = 4498
$ws.Cells.Item(14, 6).Value = 36
$ws.Cells.Item(15, 6).Value = 126
$ws.Cells.Item(19, 6).Value = 1068
$ws.Cells.Item(21, 6).Value = 3936
$ws.Cells.Item(24, 6).Value = 3308
$ws.Cells.Item(25, 6).Value = 2796
$ws.Cells.Item(26, 6).Value = 105
$ws.Cells.Item(28, 6).Value = 124
$ws.Cells.Item(30, 6).Value = 155
$ws.Cells.Item(31, 6).Value = 66
$ws.Cells.Item(35, 6).Value = 123
$ws.Cells.Item(37, 6).Value = 5085
$ws.Cells.Item(39, 6).Value = 679
$ws.Cells.Item(43, 6).Value = 952
$ws.Cells.Item(44, 6).Value = 989
$ws.Cells.Item(45, 6).Value = 385
$ws.Cells.Item(47, 6).Value = 1876
$ws.Cells.Item(49, 6).Value = 48
$ws.Cells.Item(50, 6).Value = 669
$ws.Cells.Item(51, 6).Value = 797
